$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46043,12613.5794816428,12546.4858745285,22275.86,7726.01447945948,-83.4733185838333),
    @(46044,12684.5273147281,12543.0074149232,12075.86,8597.97921998118,377.713609787684),
    @(46045,12670.357832187,11654.6660434046,12075.86,8580.35666996051,339.965113056878),
    @(46046,4859.17389324718,7754.51293041539,12075.86,7963.76503554977,151.767415248548),
    @(46047,4710.88047876496,7677.27186094577,12075.86,7895.86172170012,145.719732610245),
    @(46048,12315.9682835607,11710.6050601599,12075.86,8314.31343455133,331.210770612968),
    @(46049,12315.9682835607,12161.1379217022,12075.86,8314.31343455133,349.982973177229),
    @(46050,12315.9682835607,11996.2344983027,12075.86,8314.31343455133,343.111997202251),
    @(46051,12315.9682835607,12158.4066817302,12075.86,8314.31343455133,349.86917151173),
    @(46052,12315.9682835607,11390.0022753957,12075.86,8314.31343455133,317.852321247794),
    @(46053,4867.38022112383,7987.78371926412,12075.86,7930.25164949425,160.090640364932),
    @(46054,5046.61366744637,7713.99940170972,9743.86,7890.31053228124,244.18541391629),
    @(46055,11445.6138712783,10818.0768001126,9743.86,7806.02339765674,370.010008240389),
    @(46056,11445.6138712783,11195.4554079398,9743.86,7806.02339765674,385.734116899857)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
